$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -6
    4  = 7
    5  = -1
    7  = -4
    9  = -10
    13 = -4
    14 = -5
    17 = 9
    18 = -6
    22 = -5
    23 = -2
    24 = 3
    25 = -4
    26 = -2
    27 = -3
    28 = -2
    29 = -1
    30 = 8
    31 = 1
    33 = -3
    35 = 3
    36 = -1
    38 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
